$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update country labels that shifted position due to re-sorting by "Casos totales"
$ws.Range("A16").Value = "India"
$ws.Range("A17").Value = "Peru"
$ws.Range("A47").Value = "Sudafrica"
$ws.Range("A49").Value = "Chequia"
$ws.Range("A50").Value = "Egipto"
$ws.Range("A100").Value = "Mayotte"
$ws.Range("A101").Value = "Albania"
$ws.Range("A102").Value = "Sri Lanka"
$ws.Range("A103").Value = "Guatemala"
$ws.Range("A104").Value = "Libano"
$ws.Range("A105").Value = "Niger"
$ws.Range("A106").Value = "Costa Rica"
$ws.Range("A107").Value = "Principado de Andorra"
$ws.Range("A117").Value = "Jordania"
$ws.Range("A118").Value = "Malta"
$ws.Range("A191").Value = "Nueva Caledonia"
$ws.Range("A192").Value = "Belice"
$ws.Range("A198").Value = "Dominica"
$ws.Range("A199").Value = "Curazao"
$ws.Range("A205").Value = "Seychelles"
$ws.Range("A206").Value = "Montserrat"

# Update statistic values (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
$ws.Range("B4").Value = 1279391
$ws.Range("C4").Value = 16299
$ws.Range("E4").Value = 988909
$ws.Range("G4").Value = 1407
$ws.Range("H4").Value = 76206
$ws.Range("B9").Value = 174791
$ws.Range("C9").Value = 600
$ws.Range("D9").Value = 55027
$ws.Range("E9").Value = 93777
$ws.Range("F9").Value = 2961
$ws.Range("G9").Value = 178
$ws.Range("H9").Value = 25987
$ws.Range("B10").Value = 169015
$ws.Range("C10").Value = 853
$ws.Range("E10").Value = 21773
$ws.Range("G10").Value = 67
$ws.Range("H10").Value = 7342
$ws.Range("B16").Value = 56325
$ws.Range("C16").Value = 3338
$ws.Range("D16").Value = 16776
$ws.Range("E16").Value = 37660
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 104
$ws.Range("H16").Value = 1889
$ws.Range("B17").Value = 54817
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 17527
$ws.Range("E17").Value = 35757
$ws.Range("F17").Value = 717
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 1533
$ws.Range("B47").Value = 8232
$ws.Range("C47").Value = 424
$ws.Range("D47").Value = 3153
$ws.Range("E47").Value = 4918
$ws.Range("F47").Value = 36
$ws.Range("G47").Value = 8
$ws.Range("H47").Value = 161
$ws.Range("B48").Value = 8015
$ws.Range("C48").Value = 19
$ws.Range("E48").Value = 7766
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = 217
$ws.Range("B49").Value = 8002
$ws.Range("C49").Value = 28
$ws.Range("D49").Value = 4369
$ws.Range("E49").Value = 3364
$ws.Range("F49").Value = 52
$ws.Range("G49").Value = 7
$ws.Range("H49").Value = 269
$ws.Range("B50").Value = 7981
$ws.Range("C50").Value = 393
$ws.Range("D50").Value = 1887
$ws.Range("E50").Value = 5612
$ws.Range("F50").Value = 41
$ws.Range("G50").Value = 13
$ws.Range("H50").Value = 482
$ws.Range("B61").Value = 4199
$ws.Range("C61").Value = 265
$ws.Range("D61").Value = 2000
$ws.Range("E61").Value = 2191
$ws.Range("B80").Value = 1801
$ws.Range("C80").Value = 2
$ws.Range("D80").Value = 1755
$ws.Range("E80").Value = 36
$ws.Range("B100").Value = 854
$ws.Range("C100").Value = 115
$ws.Range("D100").Value = 352
$ws.Range("E100").Value = 492
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 10
$ws.Range("B101").Value = 842
$ws.Range("C101").Value = 10
$ws.Range("D101").Value = 605
$ws.Range("E101").Value = 206
$ws.Range("F101").Value = 7
$ws.Range("H101").Value = 31
$ws.Range("B102").Value = 805
$ws.Range("C102").Value = 8
$ws.Range("D102").Value = 232
$ws.Range("E102").Value = 564
$ws.Range("F102").Value = 1
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 9
$ws.Range("B103").Value = 798
$ws.Range("C103").Value = 35
$ws.Range("D103").Value = 86
$ws.Range("E103").Value = 691
$ws.Range("F103").Value = 5
$ws.Range("G103").Value = 2
$ws.Range("H103").Value = 21
$ws.Range("B104").Value = 784
$ws.Range("C104").Value = 34
$ws.Range("D104").Value = 220
$ws.Range("E104").Value = 539
$ws.Range("F104").Value = 43
$ws.Range("H104").Value = 25
$ws.Range("B105").Value = 770
$ws.Range("D105").Value = 561
$ws.Range("E105").Value = 171
$ws.Range("F105").Value = 0
$ws.Range("H105").Value = 38
$ws.Range("B106").Value = 761
$ws.Range("D106").Value = 428
$ws.Range("E106").Value = 327
$ws.Range("F106").Value = 5
$ws.Range("H106").Value = 6
$ws.Range("B107").Value = 751
$ws.Range("D107").Value = 521
$ws.Range("E107").Value = 184
$ws.Range("F107").Value = 15
$ws.Range("H107").Value = 46
$ws.Range("B112").Value = 648
$ws.Range("C112").Value = 31
$ws.Range("E112").Value = 626
$ws.Range("B117").Value = 494
$ws.Range("C117").Value = 21
$ws.Range("D117").Value = 381
$ws.Range("E117").Value = 104
$ws.Range("F117").Value = 5
$ws.Range("H117").Value = 9
$ws.Range("B118").Value = 486
$ws.Range("C118").Value = 2
$ws.Range("D118").Value = 413
$ws.Range("E118").Value = 68
$ws.Range("F118").Value = 0
$ws.Range("H118").Value = 5
$ws.Range("D191").Value = 18
$ws.Range("H191").Value = 0
$ws.Range("D192").Value = 16
$ws.Range("H192").Value = 2
$ws.Range("D198").Value = 14
$ws.Range("H198").Value = 0
$ws.Range("E199").Value = 1
$ws.Range("H199").Value = 1
$ws.Range("D205").Value = 8
$ws.Range("F205").Value = 0
$ws.Range("H205").Value = 0
$ws.Range("D206").Value = 7
$ws.Range("F206").Value = 1
$ws.Range("H206").Value = 1

# Update the last-updated timestamp note
$ws.Range("A1").Value = "Datos actualizados a 7 de Mayo de 2020 a las 21:04"